$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B100").Value = 6236614
$ws.Range("F100").Value = "Mineros"
$ws.Range("G100").Value = "Angostura FC"
$ws.Range("H100").Value = 1
$ws.Range("K100").Value = 2.45
$ws.Range("L100").Value = 3.3
$ws.Range("M100").Value = 2.55
$ws.Range("N100").Value = 1.8
$ws.Range("O100").Value = 3.75
$ws.Range("P100").Value = 3.6
$ws.Range("Q100").Value = -0.5
$ws.Range("R100").Value = 1.825
$ws.Range("S100").Value = 1.975
$ws.Range("U100").Value = 1.8
$ws.Range("V100").Value = 2
$ws.Range("Y100").Value = 2.6
$ws.Range("AA100").Value = 0.9750000000000001
$ws.Range("AB100").Value = 0.4
$ws.Range("AC100").Value = -0.5
$ws.Range("B101").Value = 6236257
$ws.Range("F101").Value = "CD Hermanos Colmenares"
$ws.Range("G101").Value = "Zamora"
$ws.Range("H101").Value = 0
$ws.Range("K101").Value = 2.3
$ws.Range("L101").Value = 3.2
$ws.Range("M101").Value = 2.8
$ws.Range("N101").Value = 1.666
$ws.Range("O101").Value = 3.8
$ws.Range("P101").Value = 4.2
$ws.Range("Q101").Value = -0.75
$ws.Range("R101").Value = 1.9
$ws.Range("S101").Value = 1.9
$ws.Range("U101").Value = 1.9
$ws.Range("V101").Value = 1.9
$ws.Range("Y101").Value = 3.2
$ws.Range("AA101").Value = 0.8999999999999999
$ws.Range("AB101").Value = -1
$ws.Range("AC101").Value = 0.8999999999999999
$ws.Range("B102").Value = 6236615
$ws.Range("F102").Value = "Deportivo Rayo Zuliano"
$ws.Range("G102").Value = "Academia Puerto Cabello"
$ws.Range("H102").Value = 1
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 2.375
$ws.Range("L102").Value = 3.3
$ws.Range("M102").Value = 2.625
$ws.Range("N102").Value = 2.45
$ws.Range("P102").Value = 2.55
$ws.Range("Q102").Value = 0
$ws.Range("R102").Value = 1.875
$ws.Range("S102").Value = 1.925
$ws.Range("U102").Value = 2
$ws.Range("V102").Value = 1.8
$ws.Range("W102").Value = 1.45
$ws.Range("Z102").Value = 0.875
$ws.Range("AB102").Value = -1
$ws.Range("AC102").Value = 0.8
$ws.Range("B103").Value = 6236616
$ws.Range("F103").Value = "UCV"
$ws.Range("G103").Value = "Metropolitanos FC"
$ws.Range("H103").Value = 3
$ws.Range("I103").Value = 2
$ws.Range("K103").Value = 3.3
$ws.Range("L103").Value = 3.2
$ws.Range("M103").Value = 2.05
$ws.Range("N103").Value = 2.75
$ws.Range("P103").Value = 2.3
$ws.Range("Q103").Value = 0.25
$ws.Range("R103").Value = 1.75
$ws.Range("S103").Value = 2.05
$ws.Range("U103").Value = 1.975
$ws.Range("V103").Value = 1.825
$ws.Range("W103").Value = 1.75
$ws.Range("Z103").Value = 0.75
$ws.Range("AB103").Value = 0.9750000000000001
$ws.Range("AC103").Value = -1
$ws.Range("B114").Value = 7352251
$ws.Range("F114").Value = "Caracas"
$ws.Range("G114").Value = "Academia Puerto Cabello"
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = "H"
$ws.Range("K114").Value = 2.1
$ws.Range("L114").Value = 3.2
$ws.Range("M114").Value = 3.3
$ws.Range("N114").Value = 2.15
$ws.Range("O114").Value = 3.1
$ws.Range("P114").Value = 3.2
$ws.Range("Q114").Value = -0.5
$ws.Range("R114").Value = 2.025
$ws.Range("S114").Value = 1.775
$ws.Range("T114").Value = 2.25
$ws.Range("U114").Value = 1.975
$ws.Range("V114").Value = 1.825
$ws.Range("W114").Value = 1.15
$ws.Range("X114").Value = -1
$ws.Range("Z114").Value = 1.025
$ws.Range("AA114").Value = -1
$ws.Range("AB114").Value = -1
$ws.Range("AC114").Value = 0.825
$ws.Range("B115").Value = 7352250
$ws.Range("F115").Value = "Portuguesa"
$ws.Range("G115").Value = "Deportivo Tachira"
$ws.Range("I115").Value = 1
$ws.Range("J115").Value = "D"
$ws.Range("K115").Value = 3.1
$ws.Range("L115").Value = 2.875
$ws.Range("M115").Value = 2.3
$ws.Range("N115").Value = 3
$ws.Range("O115").Value = 2.875
$ws.Range("P115").Value = 2.375
$ws.Range("Q115").Value = 0.25
$ws.Range("R115").Value = 1.725
$ws.Range("S115").Value = 2.075
$ws.Range("T115").Value = 2
$ws.Range("U115").Value = 1.825
$ws.Range("V115").Value = 1.975
$ws.Range("W115").Value = -1
$ws.Range("X115").Value = 1.875
$ws.Range("Z115").Value = 0.3625
$ws.Range("AA115").Value = -0.5
$ws.Range("AB115").Value = 0
$ws.Range("AC115").Value = -0
$ws.Range("B116").Value = 7352254
$ws.Range("F116").Value = "Academia Puerto Cabello"
$ws.Range("G116").Value = "Portuguesa"
$ws.Range("K116").Value = 2.05
$ws.Range("L116").Value = 3.4
$ws.Range("M116").Value = 3
$ws.Range("N116").Value = 1.833
$ws.Range("O116").Value = 3.5
$ws.Range("P116").Value = 3.5
$ws.Range("R116").Value = 1.65
$ws.Range("S116").Value = 2.2
$ws.Range("T116").Value = 2.25
$ws.Range("U116").Value = 1.825
$ws.Range("V116").Value = 1.975
$ws.Range("X116").Value = 2.5
$ws.Range("AA116").Value = 0.6000000000000001
$ws.Range("AB116").Value = -0.5
$ws.Range("AC116").Value = 0.4875
$ws.Range("B117").Value = 7352252
$ws.Range("F117").Value = "Deportivo Tachira"
$ws.Range("G117").Value = "Caracas"
$ws.Range("K117").Value = 2.3
$ws.Range("L117").Value = 2.875
$ws.Range("M117").Value = 3.1
$ws.Range("N117").Value = 2.25
$ws.Range("O117").Value = 2.8
$ws.Range("P117").Value = 3.25
$ws.Range("R117").Value = 1.975
$ws.Range("S117").Value = 1.825
$ws.Range("T117").Value = 2
$ws.Range("U117").Value = 1.925
$ws.Range("V117").Value = 1.875
$ws.Range("X117").Value = 1.8
$ws.Range("AA117").Value = 0.4125
$ws.Range("AB117").Value = 0
$ws.Range("AC117").Value = -0
$ws.Range("B173").Value = 7958192
$ws.Range("F173").Value = "Deportivo Tachira"
$ws.Range("G173").Value = "Monagas"
$ws.Range("H173").Value = 1
$ws.Range("J173").Value = "H"
$ws.Range("K173").Value = 1.666
$ws.Range("L173").Value = 3.4
$ws.Range("M173").Value = 4.5
$ws.Range("N173").Value = 1.95
$ws.Range("O173").Value = 3.25
$ws.Range("P173").Value = 3.5
$ws.Range("Q173").Value = -0.5
$ws.Range("R173").Value = 1.975
$ws.Range("S173").Value = 1.825
$ws.Range("W173").Value = 0.95
$ws.Range("X173").Value = -1
$ws.Range("Z173").Value = 0.9750000000000001
$ws.Range("AA173").Value = -1
$ws.Range("B174").Value = 7958193
$ws.Range("F174").Value = "Zamora"
$ws.Range("G174").Value = "Academia Puerto Cabello"
$ws.Range("H174").Value = 0
$ws.Range("J174").Value = "D"
$ws.Range("K174").Value = 3.75
$ws.Range("L174").Value = 3.3
$ws.Range("M174").Value = 1.85
$ws.Range("N174").Value = 3.1
$ws.Range("O174").Value = 3.2
$ws.Range("P174").Value = 2.1
$ws.Range("Q174").Value = 0.25
$ws.Range("R174").Value = 1.875
$ws.Range("S174").Value = 1.925
$ws.Range("W174").Value = -1
$ws.Range("X174").Value = 2.2
$ws.Range("Z174").Value = 0.4375
$ws.Range("AA174").Value = -0.5
$ws.Range("N183").Value = 2
$ws.Range("O183").Value = 3.2
$ws.Range("P183").Value = 3.75
$ws.Range("Q183").Value = -0.5
$ws.Range("R183").Value = 2.025
$ws.Range("S183").Value = 1.775
$ws.Range("U183").Value = 1.975
$ws.Range("V183").Value = 1.825
$ws.Range("R185").Value = 2.025
$ws.Range("S185").Value = 1.775
